$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("kole_warren_arjun191@protonmail.com", "Arjun Warren Kole ", "kole_warren_arjun191", "TqoVitqiYze"),
    @("levi_elisha_tony387@gmail.com", "Tony Elisha Levi ", "levi_elisha_tony387", "7GYKCz`$dl"),
    @("julian_gunner_cooper673@yahoo.com", "Cooper Gunner Julian ", "julian_gunner_cooper673", "DXGkQ*!yRo"),
    @("sean_dash_kasen279@inbox.com", "Kasen Dash Sean ", "sean_dash_kasen279", "aitmjNs)QM0"),
    @("isaias_elijah_moises979@yahoo.com", "Moises Elijah Isaias ", "isaias_elijah_moises979", "5SVq@wKz2Al"),
    @("eli_terry_graham532@hotmail.com", "Graham Terry Eli ", "eli_terry_graham532", "kdJp2u"),
    @("matthew_gibson_nolan260@inbox.com", "Nolan Gibson Matthew ", "matthew_gibson_nolan260", "UCvIH*LH6")
)

$startRow = 27
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}

$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("C35").Select()
